$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A4").Value = -20.726
$ws.Range("A6").Value = -22.291
$ws.Range("E6").Value = 16.451
$ws.Range("A7").Value = -21.304
$ws.Range("B7").Value = 5.872
$ws.Range("E8").Value = 16.344
$ws.Range("B12").Value = 4.786
$ws.Range("D12").Value = -7.393000000000001
$ws.Range("C13").Value = -13.072
$ws.Range("C14").Value = -11.818
$ws.Range("B15").Value = 5.086999999999999
$ws.Range("A16").Value = -21.632
$ws.Range("C16").Value = -13.107
$ws.Range("E18").Value = 16.661
$ws.Range("C19").Value = -12.505
$ws.Range("E19").Value = 16.891
$ws.Range("A20").Value = -21.433
$ws.Range("B20").Value = 6.25
$ws.Range("E20").Value = 16.408
$ws.Range("B21").Value = 8.074
$ws.Range("E21").Value = 16.78
$ws.Range("B22").Value = 7.894
$ws.Range("C22").Value = -12.787
$ws.Range("D22").Value = -8.120999999999999
$ws.Range("B23").Value = 7.369000000000002
$ws.Range("E23").Value = 16.304
$ws.Range("E24").Value = 16.731
$ws.Range("A28").Value = -21.687
$ws.Range("A29").Value = -21.286
$ws.Range("B29").Value = 6.180999999999999
$ws.Range("D29").Value = -7.128
$ws.Range("A32").Value = -21.486
$ws.Range("B34").Value = 7.821
$ws.Range("D34").Value = -7.929
$ws.Range("E35").Value = 16.294
$ws.Range("C36").Value = -12.683
$ws.Range("E37").Value = 16.509
$ws.Range("E39").Value = 17.071
$ws.Range("A40").Value = -20.486
$ws.Range("E41").Value = 16.597
$ws.Range("B42").Value = 7.582999999999998
$ws.Range("B43").Value = 5.907000000000001
$ws.Range("D43").Value = -8.145000000000001
$ws.Range("B44").Value = 5.171
$ws.Range("B45").Value = 5.048
$ws.Range("A46").Value = -21.519
$ws.Range("B46").Value = 6.201000000000001
$ws.Range("C46").Value = -13.323
$ws.Range("E46").Value = 16.723
$ws.Range("D48").Value = -7.797
$ws.Range("B50").Value = 5.167999999999999
$ws.Range("C50").Value = -13.097
$ws.Range("A51").Value = -21.292
$ws.Range("B51").Value = 6.7
$ws.Range("A52").Value = -21.561
$ws.Range("A57").Value = -22.028
$ws.Range("E58").Value = 16.489
$ws.Range("A59").Value = -22.248
$ws.Range("D60").Value = -8.178999999999998
$ws.Range("E60").Value = 16.591
$ws.Range("A62").Value = -21.937
$ws.Range("A66").Value = -21.54
$ws.Range("B66").Value = 5.809000000000001
$ws.Range("B67").Value = 5.249
$ws.Range("D68").Value = -6.900999999999999
$ws.Range("D70").Value = -7.103
$ws.Range("A73").Value = -19.95
$ws.Range("D73").Value = -8.345000000000001
$ws.Range("E73").Value = 16.487
$ws.Range("A74").Value = -21.061
$ws.Range("E76").Value = 16.457
$ws.Range("B79").Value = 5.529999999999999
$ws.Range("B84").Value = 5.798
$ws.Range("E85").Value = 16.83
$ws.Range("D87").Value = -8.475
$ws.Range("A92").Value = -20.898
$ws.Range("B92").Value = 6.153999999999999
$ws.Range("D92").Value = -6.497
$ws.Range("C95").Value = -11.719
$ws.Range("B97").Value = 6.188999999999999
$ws.Range("C97").Value = -13.247
$ws.Range("E98").Value = 16.406
$ws.Range("A100").Value = -21.526
$ws.Range("D101").Value = -7.890000000000001
